$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing data before rearranging (rows 2-6, cols A:B)
$names = @{}
$emails = @{}
for ($r = 2; $r -le 6; $r++) {
    $names[$r] = $ws.Cells.Item($r, 1).Value()
    $emails[$r] = $ws.Cells.Item($r, 2).Value()
}

# Row 6 (Thiago) moves up to become row 2 - the primary notification contact
$ws.Cells.Item(2, 1).Value = $names[6]
$ws.Cells.Item(2, 2).Value = $emails[6]

# Clear out the old rows 3-6 in columns A:B (data relocated to N:O)
$ws.Range("A3:B6").ClearContents()

# Remaining team members (originally rows 2-5: Amanda, Bruna, Luciana, Nathalia)
# are relocated to columns N:O, starting at row 3
$ws.Cells.Item(3, 14).Value = $names[2]
$ws.Cells.Item(3, 15).Value = $emails[2]

$ws.Cells.Item(4, 14).Value = $names[3]
$ws.Cells.Item(4, 15).Value = $emails[3]

$ws.Cells.Item(5, 14).Value = $names[4]
$ws.Cells.Item(5, 15).Value = $emails[4]

$ws.Cells.Item(6, 14).Value = $names[5]
$ws.Cells.Item(6, 15).Value = $emails[5]

# Update the active selection to match the new layout
$ws.Range("B9").Select()
